# Updated main GSC export data: the 2025-10-19 row (which had no
# "No video indexed" / "Video indexed" counts yet) is removed from the
# "Chart" sheet, and every subsequent day's row shifts up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows("2").Delete()
